$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 146..148, pushing existing rows 146-176 down to 149-179.
$ws.Range("A146:R148").EntireRow.Insert()

# Row 146 (new): Banquete, $/bandeja 10 kilos
$ws.Cells.Item(146, 1).Value = 9
$ws.Cells.Item(146, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(146, 3).Value = "Metropolitana"
$ws.Cells.Item(146, 4).Value = 45211
$ws.Cells.Item(146, 5).Value = 13
$ws.Cells.Item(146, 6).Value = 300000000
$ws.Cells.Item(146, 7).Value = "Espárragos"
$ws.Cells.Item(146, 8).Value = "Sin especificar"
$ws.Cells.Item(146, 9).Value = "Banquete"
$ws.Cells.Item(146, 10).Value = 52
$ws.Cells.Item(146, 11).Value = 16000
$ws.Cells.Item(146, 12).Value = 16000
$ws.Cells.Item(146, 13).Value = 16000
$ws.Cells.Item(146, 14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(146, 15).Value = "Provincia de Linares"
$ws.Cells.Item(146, 16).Value = 1600
$ws.Cells.Item(146, 17).Value = 10
$ws.Cells.Item(146, 18).Value = "Hortaliza"

# Row 147 (new): Primera, $/bandeja 10 kilos
$ws.Cells.Item(147, 1).Value = 9
$ws.Cells.Item(147, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(147, 3).Value = "Metropolitana"
$ws.Cells.Item(147, 4).Value = 45211
$ws.Cells.Item(147, 5).Value = 13
$ws.Cells.Item(147, 6).Value = 300000000
$ws.Cells.Item(147, 7).Value = "Espárragos"
$ws.Cells.Item(147, 8).Value = "Sin especificar"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 160
$ws.Cells.Item(147, 11).Value = 14000
$ws.Cells.Item(147, 12).Value = 14000
$ws.Cells.Item(147, 13).Value = 14000
$ws.Cells.Item(147, 14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(147, 15).Value = "Provincia de Linares"
$ws.Cells.Item(147, 16).Value = 1400
$ws.Cells.Item(147, 17).Value = 10
$ws.Cells.Item(147, 18).Value = "Hortaliza"

# Row 148 (new): Segunda, $/bandeja 10 kilos
$ws.Cells.Item(148, 1).Value = 9
$ws.Cells.Item(148, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(148, 3).Value = "Metropolitana"
$ws.Cells.Item(148, 4).Value = 45211
$ws.Cells.Item(148, 5).Value = 13
$ws.Cells.Item(148, 6).Value = 300000000
$ws.Cells.Item(148, 7).Value = "Espárragos"
$ws.Cells.Item(148, 8).Value = "Sin especificar"
$ws.Cells.Item(148, 9).Value = "Segunda"
$ws.Cells.Item(148, 10).Value = 70
$ws.Cells.Item(148, 11).Value = 12000
$ws.Cells.Item(148, 12).Value = 12000
$ws.Cells.Item(148, 13).Value = 12000
$ws.Cells.Item(148, 14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(148, 15).Value = "Provincia de Linares"
$ws.Cells.Item(148, 16).Value = 1200
$ws.Cells.Item(148, 17).Value = 10
$ws.Cells.Item(148, 18).Value = "Hortaliza"
